$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 112065913
$ws.Range("B2").Value = 90799
$ws.Range("E2").Value = 1968
$ws.Range("F2").Value = "Grantaggsvamp"
$ws.Range("G2").Value = "Bankera violascens"
$ws.Range("H2").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q2").Value = 469597
$ws.Range("R2").Value = 7039829

# --- Row 3 ---
$ws.Range("B3").Value = 89072

# --- Row 4 ---
$ws.Range("A4").Value = 112065865
$ws.Range("B4").Value = 90806
$ws.Range("E4").Value = 4361
$ws.Range("F4").Value = "Orange taggsvamp"
$ws.Range("G4").Value = "Hydnellum aurantiacum"
$ws.Range("H4").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q4").Value = 469609
$ws.Range("R4").Value = 7039805

# --- Row 5 ---
$ws.Range("A5").Value = 112370019
$ws.Range("B5").Value = 56430
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = "Tretåig hackspett"
$ws.Range("G5").Value = "Picoides tridactylus"
$ws.Range("H5").Value = "(Linnaeus, 1758)"
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("Q5").Value = 469645
$ws.Range("R5").Value = 7039915
$ws.Range("AC5").Value = "ringhack äldre"

# --- Row 6 ---
$ws.Range("A6").Value = 112370026
$ws.Range("B6").Value = 90806
$ws.Range("E6").Value = 4361
$ws.Range("F6").Value = "Orange taggsvamp"
$ws.Range("G6").Value = "Hydnellum aurantiacum"
$ws.Range("H6").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("K6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("Q6").Value = 469718
$ws.Range("R6").Value = 7039994
$ws.Range("AC6").ClearContents()
